$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Use case 1 & 2: cell B6 previously held the text "user15" (a shared string).
# It is now changed to the numeric value 1.
$ws.Range("B6").Value = 1

# Update the selection to reflect the edited cell.
$ws.Range("B6").Select()
